# Add a new "2022-Q3" sheet right after "总计" and insert its summary row
# into the "总计" sheet, shifting all other quarters down by one.

$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert a new row for 2022-Q3
#    right after the header row, pushing the existing data rows down.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Remember the formatting used for the "index" column (A) and the header
# row so we can reapply it to the new worksheet further down.
$indexStyleCell = $summary.Cells.Item(2, 1)   # A2, style used for the row index column
$headerStyleCell = $summary.Cells.Item(1, 2)  # B1, bold/centered header style

# Shift existing data rows 2..8 down to 3..9 (copy values+format, bottom-up
# so we never overwrite a row before it has been copied).
for ($r = 8; $r -ge 2; $r--) {
    $srcRow = $summary.Range("A" + $r + ":D" + $r)
    $dstRow = $summary.Range("A" + ($r + 1) + ":D" + ($r + 1))
    $srcRow.Copy($dstRow)
}

# Fill in the new row 2 with the 2022-Q3 figures.
$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 9
$summary.Cells.Item(2, 4).Value = 3.72

# ---------------------------------------------------------------------
# 2. Insert a brand-new worksheet named "2022-Q3" right after "总计"
#    holding the fund-level detail for that quarter.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

# Helper to write a value as text even when it looks numeric (keeps the
# same cell "type" the source workbook uses for these text-formatted
# columns), then restore the default (no explicit) number format.
function Set-TextValue($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Header row (bold/centered like the other quarterly sheets).
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Count; $c++) {
    $q3.Cells.Item(1, 2 + $c).Value = $headers[$c]
}
$headerStyleCell.Copy()
$q3.Range("B1:H1").PasteSpecial($xlPasteFormats)

# Data rows.
$rows = @(
    @("001822", "华商智能生活灵活配置混合A", "33.45", "87.34", "4.59", "1.5354", 6),
    @("001933", "华商新兴活力灵活配置混合", "18.34", "87.96", "4.55", "0.8345", 5),
    @("015385", "华商智能生活灵活配置混合C", "11.97", "87.34", "4.59", "0.5494", 6),
    @("013886", "华商新能源汽车混合A", "9.72", "89.01", "4.70", "0.4568", 8),
    @("013887", "华商新能源汽车混合C", "4.06", "89.01", "4.70", "0.1908", 8),
    @("014350", "华商卓越成长一年持有混合A", "3.14", "86.88", "4.54", "0.1426", 8),
    @("014351", "华商卓越成长一年持有混合C", "0.10", "86.88", "4.54", "0.0045", 8),
    @("015466", "太平中证1000指数增强A", "0.37", "92.23", "1.05", "0.0039", 6),
    @("015467", "太平中证1000指数增强C", "0.02", "92.23", "1.05", "0.0002", 6)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $q3.Cells.Item($r, 1).Value = $i

    Set-TextValue $q3.Cells.Item($r, 2) $data[0]
    Set-TextValue $q3.Cells.Item($r, 3) $data[1]
    Set-TextValue $q3.Cells.Item($r, 4) $data[2]
    Set-TextValue $q3.Cells.Item($r, 5) $data[3]
    Set-TextValue $q3.Cells.Item($r, 6) $data[4]
    Set-TextValue $q3.Cells.Item($r, 7) $data[5]

    $q3.Cells.Item($r, 8).Value = $data[6]
}

# Apply the "index column" style (bold/centered with border) to column A
# of every data row, matching the other quarterly sheets.
$indexStyleCell.Copy()
$q3.Range("A2:A10").PasteSpecial($xlPasteFormats)

$q3.Range("A1").Select()
